
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Split the tracked "io/html-portfolio" insertion so that only "io/"
#    remains part of the original tracked insertion (w:ins id=2), and the
#    rest of the old link text ("html-portfolio") is replaced by a brand new,
#    *untracked* run containing "HTML-Portfolio". A plain Find/Replace over
#    the tail of that tracked run produces exactly this split.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("html-portfolio", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "HTML-Portfolio", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Relocate the "_GoBack" bookmark from the very end of the document to the
#    end of the paragraph that now reads ".../io/HTML-Portfolio" (i.e. right
#    after the newly typed run, before the paragraph mark).
# ---------------------------------------------------------------------------

# Find that paragraph dynamically (don't hard-code its index).
$linkPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*HTML-Portfolio*") {
        $linkPara = $p
    }
}

# Position right before the paragraph mark of that paragraph.
$boundary = $linkPara.Range.End - 1

# A bookmark collapsed exactly on a paragraph-end boundary cannot be added
# directly and reliably, so nudge it into place: insert a sentinel
# character, drop the (now non-boundary) collapsed bookmark in front of it,
# then delete the sentinel again.
$sentinelSpot = $d.Range($boundary, $boundary)
$sentinelSpot.InsertAfter("Z")

$bmSpot = $d.Range($boundary, $boundary)
$d.Bookmarks.Add("_GoBack", $bmSpot) | Out-Null

$sentinelRange = $d.Range($boundary, $boundary + 1)
$sentinelRange.Delete() | Out-Null

Write-Host "Updated GitHub Pages link and relocated the _GoBack bookmark."
